$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to reflect the latest scrape.
# NumberFormat is forced to text ("@") before assignment so that numeric-looking
# strings (e.g. "0.647") are stored as text, matching the source data, then the
# cell style is reset back to "Normal" so no stray formatting is introduced.

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '43.012.66'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +0.50%  '
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.301.90'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +2.25%  '
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.10%  '
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '251.81'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +0.31%  '
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.647'
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +2.16%  '
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '75.27'
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +7.08%  '
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.652'
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -4.25%  '
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '39.31'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +0.26%  '
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0986'
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +1.32%  '
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '7.58'
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -0.13%  '
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.106'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +1.49%  '
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '2.647.08'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +2.17%  '
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '15.17'
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +2.08%  '
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -0.82%  '
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '2.301.90'
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +1.99%  '
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '42.918.24'
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +0.46%  '
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.0000100'
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +2.10%  '
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.28'
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -0.21%  '
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '72.47'
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -0.73%  '
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '237.57'
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +0.82%  '
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '2.19'
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +6.14%  '
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '3.88'
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -0.94%  '
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '11.44'
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  -2.43%  '
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +0.13%  '
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '2.39'
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -2.41%  '
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.14'
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -0.44%  '
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '167.95'
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +0.33%  '
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '21.07'
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +0.35%  '
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.0864'
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +10.07%  '
$c.Style = 'Normal'
$c = $ws.Range('B32')
$c.NumberFormat = '@'
$c.Value = 'InternetComputer(DFINITY)'
$c.Style = 'Normal'
$c = $ws.Range('C32')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '6.29'
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -1.04%  '
$c.Style = 'Normal'
$c = $ws.Range('B33')
$c.NumberFormat = '@'
$c.Value = 'Kaspa'
$c.Style = 'Normal'
$c = $ws.Range('C33')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.129'
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +1.02%  '
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '31.33'
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +0.56%  '
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +1.19%  '
$c.Style = 'Normal'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '4.67'
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +8.27%  '
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '4.82'
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +2.37%  '
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.0306'
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -5.47%  '
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '13.77'
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +10.81%  '
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.31'
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +0.28%  '
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '6.01'
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +3.05%  '
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.213'
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +5.55%  '
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '9.12'
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +1.58%  '
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '61.42'
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -2.84%  '
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -0.99%  '
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '105.88'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +11.74%  '
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -2.02%  '
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -0.23%  '
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.18'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -1.31%  '
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '4.26'
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -1.18%  '
$c.Style = 'Normal'
